# edit.ps1 -- Applies the "Uebersicht ueber Testspiele aktualisiert" update:
#   - Row 22 (Marcus Riemer): "Zu Verbessern" text extended with the
#     guard-turning note and a clarification about the night level.
#   - Row 23 (Arne Kaleck): unchanged.
#   - Row 24 (NEW, Jan Helke): long playtest feedback text.
#   - Row 25 (NEW, Lars' Bruder): name only, no feedback text yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 22: Marcus Riemer's "Zu Verbessern" (column F) gets rewritten.
# ---------------------------------------------------------------------
$ws.Range("F22").Value = "Drehen der Wachen verbessern (erst stoppen, dann drehen, dann weiterlaufen, etc.); Erweiterungsvorschläge: mehr Level, bessers Interface, Schalter für Sound (anstatt Stein), Nachtlevel mit eingeschränktem Sichtradius des Spielers, Beleuchtung"

# ---------------------------------------------------------------------
# Row 23 (Arne Kaleck) is left as-is.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# New row 24: Jan Helke. Version/Testmethoden copied from the row above.
# ---------------------------------------------------------------------
$ws.Range("B24").Value = $ws.Range("B23").Value2
$ws.Range("C24").Value = $ws.Range("C23").Value2
$ws.Range("D24").Value = "Jan Helke"

$spielideeLines = @(
    "Spielidee: Ich finde, das Spielprinzip „Vermeide die Sichtlinie der Gegner“ ist jetzt nicht sonderlich originell, da gab es in der Vergangenheit viele Titel, die entweder exakt das als Spielziel hatten oder bei denen diese Taktik mindestens hilfreich ist. Namen, die mir da einfallen sind „Splinter Cell“ und „Commandos – Hinter feindlichen Linien“. Aber: Ich ma das Spielprinzip, deshalb hatte ich Spaß beim Test.",
    "Hintergrund: Ich würde mir wünsche, einen minimalen Hintergrund zu haben. Irgendein einleitender Text „Du bist Spion und willst die DA ausspionieren. Pass auf, dass du dabei nicht erwischt wirst.“ Einfach, damit man nicht alleine mit dem Start-Button da sitzt und nicht weiß, was man machen muss.",
    "Schwierigkeit: Generell fand ich die Level alle recht einfach. Die Wachen wirken ziemlich statisch.",
    "Idee: Kann man den Wachen ein generelles Pattern geben, dass sie (egal, ob sie stillstehen oder sich bewegen) immer „den Kopf“ um 5 – 10 Grad in beide Richtungen drehen. Gerne auch in einem immer zufälligen Winkel. So wie jeder Mensch ja mal mehr mal weniger nach links oder rechts guckt. Ich könnte mir vorstellen, dass das als generische Funktion aller Wachen recht einfach einzubauen sein müsste.",
    "Level 1: Die obere Wache sollte sich beim Wenden in den Raum hinein drehen. Es ist unwahrscheinlich, dass sie beim Drehen immer die Wand angucken würde.",
    "Level 2: Ich finde die obere Wache ziemlich dämlich. Die sollte sich auf jeden Fall beständig langsam von links nach rechts und zurück drehen. Kann man die Disquise-Option ausblenden? Die braucht man hier noch nicht und sie ist nicht erklärt.",
    "Level 3: Auch hier würde ich die Disquise Option ausblenden. Ich würde ansonsten auch die Benutzung der Shift Taste ändern. Und zwar würde ich das „Wache durch die Gegend schleifen“ nur machen, wenn die Taste gedrückt gehalten wird. Das fühlt sich einfach haptisch besser an, dass ich die Taste halten muss, während ich die Wache ziehe. Wie beim Laufen, da muss ich ja auch Halten. Und sobald ich die Taste loslasse, lass ich die Wache exakt da liegen, wo sie liegt.",
    "Level 4: In dem Level ist mir das mit den stur nach vorne schauenden Wachen besonders aufgefallen. Es würde lebendiger wirken, wenn sie ein bisschen schwenken würden. Das hätte zwar keine Auswirkungen auf die Schwierigkeit des Levels, wirkt aber aktiver. Ansonsten habe ich festgestellt, dass sich der schwarze Ring um den Chrakter anders verhält. Im normalen Modus ist der so etwas wie eine Pufferzone. Wenn der schwarze Ring in den Sichtbereich einer Wache gerät, passiert erst mal gar nichts. Der Alarm geht erst los, wenn die Wache den inneren Ring „sieht“. Wenn der Charakter getarnt ist, werden die Wachen getriggert in dem Moment, wenn sie den schwarzen Ring sehen. Das würde ich vereinheitlichen und den schwarzen Ring bei beiden Varianten als Pufferzone nutzen.",
    "Level 5: Im Text steht zwar, dass man die Wurfweite ändern kann, je nachdem, wie lange man die Taste hält. Ich würde aber noch einen Hinweis auf die farbliche Veränderung des Steins je nach Wurfkraft einbauen.",
    "Level 6: Gefällt mir gut. Es kommen unterschiedliche Elemente zum Einsatz und es gibt verschiedene Arten, das Level zu lösen. Ich würde der oberen Wache einen etwas komplizierteren Weg geben. Und zwar dass sie sich an der rechten unteren Ecke ihres Weges umdreht und dann den Weg in die entgegengesetzte Richtung geht. Das zeigt auch, dass die Wachen etwas mehr können als nur im Krei zu laufen."
)
$ws.Range("F24").Value = [string]::Join("`r`n", $spielideeLines)

# ---------------------------------------------------------------------
# New row 25: Lars' Bruder (name only, no write-up yet).
# ---------------------------------------------------------------------
$ws.Range("B25").Value = $ws.Range("B23").Value2
$ws.Range("C25").Value = $ws.Range("C23").Value2
$ws.Range("D25").Value = "Lars' Bruder"

# ---------------------------------------------------------------------
# Row heights: the two touched/added rows take on their wrapped-text
# height; the untouched rows above settle back to their recalculated
# auto-fit heights now that the sheet has new content below them.
# ---------------------------------------------------------------------
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 75
$ws.Rows.Item(17).RowHeight = 150
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(24).RowHeight = 409.5

# ---------------------------------------------------------------------
# Move the selection/view down to the newly added row, mirroring the
# author's cursor position after entering the new feedback.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 1
$ws.Range("F24").Select()
